{"js": "// Map of old division expressions to new ones, as produced by the commit.\nconst replacements = [\n  [\"19\u00f77=\", \"27\u00f73=\"],\n  [\"54\u00f72=\", \"93\u00f75=\"],\n  [\"60\u00f79=\", \"72\u00f78=\"],\n  [\"89\u00f77=\", \"17\u00f76=\"],\n  [\"66\u00f72=\", \"49\u00f78=\"],\n  [\"58\u00f74=\", \"77\u00f76=\"],\n  [\"18\u00f77=\", \"71\u00f75=\"],\n  [\"74\u00f72=\", \"69\u00f73=\"],\n  [\"86\u00f76=\", \"39\u00f78=\"],\n  [\"72\u00f75=\", \"18\u00f78=\"],\n  [\"80\u00f78=\", \"33\u00f74=\"],\n  [\"98\u00f78=\", \"69\u00f76=\"],\n  [\"59\u00f79=\", \"58\u00f72=\"],\n  [\"13\u00f74=\", \"79\u00f78=\"],\n  [\"94\u00f73=\", \"77\u00f72=\"],\n  [\"21\u00f75=\", \"65\u00f74=\"],\n  [\"51\u00f77=\", \"30\u00f74=\"],\n  [\"27\u00f77=\", \"65\u00f77=\"],\n  [\"99\u00f77=\", \"81\u00f78=\"],\n  [\"42\u00f76=\", \"65\u00f73=\"],\n  [\"87\u00f73=\", \"99\u00f75=\"],\n  [\"61\u00f76=\", \"83\u00f74=\"],\n  [\"49\u00f73=\", \"17\u00f79=\"],\n  [\"53\u00f78=\", \"26\u00f75=\"],\n  [\"29\u00f75=\", \"64\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Map of old division expressions to new ones, as produced by the commit.\n$replacements = @(\n    @{old=\"19\u00f77=\"; new=\"27\u00f73=\"},\n    @{old=\"54\u00f72=\"; new=\"93\u00f75=\"},\n    @{old=\"60\u00f79=\"; new=\"72\u00f78=\"},\n    @{old=\"89\u00f77=\"; new=\"17\u00f76=\"},\n    @{old=\"66\u00f72=\"; new=\"49\u00f78=\"},\n    @{old=\"58\u00f74=\"; new=\"77\u00f76=\"},\n    @{old=\"18\u00f77=\"; new=\"71\u00f75=\"},\n    @{old=\"74\u00f72=\"; new=\"69\u00f73=\"},\n    @{old=\"86\u00f76=\"; new=\"39\u00f78=\"},\n    @{old=\"72\u00f75=\"; new=\"18\u00f78=\"},\n    @{old=\"80\u00f78=\"; new=\"33\u00f74=\"},\n    @{old=\"98\u00f78=\"; new=\"69\u00f76=\"},\n    @{old=\"59\u00f79=\"; new=\"58\u00f72=\"},\n    @{old=\"13\u00f74=\"; new=\"79\u00f78=\"},\n    @{old=\"94\u00f73=\"; new=\"77\u00f72=\"},\n    @{old=\"21\u00f75=\"; new=\"65\u00f74=\"},\n    @{old=\"51\u00f77=\"; new=\"30\u00f74=\"},\n    @{old=\"27\u00f77=\"; new=\"65\u00f77=\"},\n    @{old=\"99\u00f77=\"; new=\"81\u00f78=\"},\n    @{old=\"42\u00f76=\"; new=\"65\u00f73=\"},\n    @{old=\"87\u00f73=\"; new=\"99\u00f75=\"},\n    @{old=\"61\u00f76=\"; new=\"83\u00f74=\"},\n    @{old=\"49\u00f73=\"; new=\"17\u00f79=\"},\n    @{old=\"53\u00f78=\"; new=\"26\u00f75=\"},\n    @{old=\"29\u00f75=\"; new=\"64\u00f79=\"}\n)\n\n$d = $word.ActiveDocument\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($r.old, $false, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)\n}\n"}
